$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "296.42"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "3.85%"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "41.42"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "3.57%"
$r.Style = "Normal"
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "5.042"
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "-0.02%"
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.07469"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "2.79%"
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "4.357"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "1.08%"
$r.Style = "Normal"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.581"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "4.98%"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.9278"
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "1.35%"
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.1185"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "-1.15%"
$r.Style = "Normal"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.1828"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "7.40%"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.08872"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "3.67%"
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.04192"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "0.45%"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.1052"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "-0.02%"
$r.Style = "Normal"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.001276"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "2.16%"
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.005976"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "0.03%"
$r.Style = "Normal"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.003875"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "2.36%"
$r.Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "3.346"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "-1.67%"
$r.Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.3310"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "1.47%"
$r.Style = "Normal"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.882"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "0.23%"
$r.Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.1408"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "3.89%"
$r.Style = "Normal"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.2967"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "2.66%"
$r.Style = "Normal"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.04032"
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "4.69%"
$r.Style = "Normal"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.001264"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "-0.56%"
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.0001231"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "-4.15%"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.0003724"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "-0.33%"
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.02392"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "4.58%"
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.05210"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "6.10%"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.006685"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "-4.66%"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.007794"
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "1.23%"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.1322"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "4.61%"
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.007376"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "-0.26%"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.007193"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "3.82%"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.3214"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "4.22%"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.00006229"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "-2.96%"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.00000000750"
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "-0.32%"
$r.Style = "Normal"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.004202"
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "0.01%"
$r.Style = "Normal"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.04601"
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "-81.73%"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.00002101"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "-0.32%"
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0002001"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "-0.32%"
$r.Style = "Normal"
